$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.671.83"
$ws.Range("E2").Value = "  +3.73%  "

$ws.Range("D3").Value = "3.497.68"
$ws.Range("E3").Value = "  +2.06%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.24"
$ws.Range("E5").Value = "  +3.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.85"
$ws.Range("E6").Value = "  +5.00%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").Value = "3.495.58"
$ws.Range("E8").Value = "  +1.95%  "

$ws.Range("E9").Value = "  +4.35%  "

$ws.Range("E10").Value = "  +0.91%  "

$ws.Range("E11").Value = "  +4.01%  "

$ws.Range("E12").Value = "  +3.17%  "

$ws.Range("D13").Value = "4.104.84"
$ws.Range("E13").Value = "  +2.18%  "

$ws.Range("E14").Value = "  -0.52%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.25"
$ws.Range("E15").Value = "  +3.72%  "

$ws.Range("E16").Value = "  +1.41%  "

$ws.Range("D17").Value = "66.665.99"
$ws.Range("E17").Value = "  +3.60%  "

$ws.Range("D18").Value = "3.517.69"
$ws.Range("E18").Value = "  +2.38%  "

$ws.Range("E19").Value = "  +3.46%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.07"
$ws.Range("E20").Value = "  +3.12%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "389.72"
$ws.Range("E21").Value = "  +2.59%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.98"
$ws.Range("E22").Value = "  +1.11%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.09"
$ws.Range("E23").Value = "  +2.02%  "

$ws.Range("E24").Value = "  -0.12%  "

$ws.Range("E25").Value = "  +2.72%  "

$ws.Range("E26").Value = "  +4.29%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.40"
$ws.Range("E27").Value = "  +7.45%  "

$ws.Range("E28").Value = "  +2.12%  "

$ws.Range("E29").Value = "  +0.33%  "

$ws.Range("E30").Value = "  +4.14%  "

$ws.Range("E31").Value = "  +5.28%  "

$ws.Range("E32").Value = "  +2.51%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.58"
$ws.Range("E33").Value = "  +2.37%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.46"
$ws.Range("E34").Value = "  +4.90%  "

$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("E36").Value = "  +6.13%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "162.66"
$ws.Range("E37").Value = "  +2.25%  "

$ws.Range("E38").Value = "  +2.55%  "

$ws.Range("E39").Value = "  +3.15%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.72"
$ws.Range("E40").Value = "  +5.57%  "

$ws.Range("E41").Value = "  +2.10%  "

$ws.Range("D44").Value = "2.825.52"
$ws.Range("E44").Value = "  +0.31%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "26.66"
$ws.Range("E45").Value = "  +0.54%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "43.04"
$ws.Range("E46").Value = "  -0.39%  "

$ws.Range("E47").Value = "  +2.23%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "355.20"
$ws.Range("E48").Value = "  +2.38%  "

$ws.Range("E49").Value = "  +4.41%  "

$ws.Range("E50").Value = "  +2.56%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.74"
$ws.Range("E51").Value = "  +11.35%  "

# Row 42: EnergySwap -> RenderToken (full row update)
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.75"
$ws.Range("E42").Value = "  +2.66%  "

# Row 43: RenderToken -> EnergySwap (full row update)
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.45"
$ws.Range("E43").Value = "  +2.29%  "
